# Auto-generated edit script applying scheduled-runner updates to Sheets/Atomos_Profits.xlsx
# Updates currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 2759.923
$ws.Range("I62").Value = 2538.0908
$ws.Range("J62").Value = 3980
$ws.Range("K62").Value = 2538.0908
$ws.Range("L62").Value = 3980
$ws.Range("M62").Value = -1914.0908
$ws.Range("N62").Value = -5228

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 2759.923
$ws.Range("I65").Value = 2538.0908
$ws.Range("J65").Value = 3980
$ws.Range("K65").Value = 12690.454
$ws.Range("L65").Value = 19900
$ws.Range("M65").Value = -9570.454
$ws.Range("N65").Value = -26140

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3281376
$ws.Range("J132").Value = 3992.7144
$ws.Range("L132").Value = 11978.1432
$ws.Range("N132").Value = -17038.1432

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2360.7327
$ws.Range("I138").Value = 819.2406999999999
$ws.Range("J138").Value = 4962
$ws.Range("K138").Value = 2457.7221
$ws.Range("L138").Value = 14886
$ws.Range("M138").Value = 2682.2779
$ws.Range("N138").Value = -25166

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 218381.11
$ws.Range("I141").Value = 1159.3077
$ws.Range("J141").Value = 1159675.5
$ws.Range("K141").Value = 3477.9231
$ws.Range("L141").Value = 3479026.5
$ws.Range("M141").Value = 1702.0769
$ws.Range("N141").Value = -3489386.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 1378.23
$ws.Range("I32").Value = 1334.4362
$ws.Range("J32").Value = 2064.3333
$ws.Range("K32").Value = 1334.4362
$ws.Range("L32").Value = 2064.3333
$ws.Range("M32").Value = -1047.4362
$ws.Range("N32").Value = -2638.3333

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 1628.2325
$ws.Range("I61").Value = 828.5806
$ws.Range("J61").Value = 3694
$ws.Range("K61").Value = 828.5806
$ws.Range("L61").Value = 3694
$ws.Range("M61").Value = -616.5806
$ws.Range("N61").Value = -4118

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1894.1964
$ws.Range("I132").Value = 1610.093
$ws.Range("J132").Value = 2833.923
$ws.Range("K132").Value = 4830.279
$ws.Range("L132").Value = 8501.769
$ws.Range("M132").Value = -2300.279
$ws.Range("N132").Value = -13561.769

# Row 134: Brace for More Vambraces
$ws.Range("H134").Value = 34625
$ws.Range("J134").Value = 34625
$ws.Range("L134").Value = 34625
$ws.Range("N134").Value = -44765

# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 33998.5
$ws.Range("J135").Value = 33998.5
$ws.Range("L135").Value = 33998.5
$ws.Range("N135").Value = -44138.5

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 1628.2325
$ws.Range("I136").Value = 828.5806
$ws.Range("J136").Value = 3694
$ws.Range("K136").Value = 2485.7418
$ws.Range("L136").Value = 11082
$ws.Range("M136").Value = 64.25820000000022
$ws.Range("N136").Value = -16182

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 17686.188
$ws.Range("I86").Value = 1240.0454
$ws.Range("K86").Value = 1240.0454
$ws.Range("M86").Value = -117.0454

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 17686.188
$ws.Range("I89").Value = 1240.0454
$ws.Range("K89").Value = 6200.227
$ws.Range("M89").Value = -584.2269999999999

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 2576.4736
$ws.Range("I99").Value = 1976.3334
$ws.Range("J99").Value = 3605.2856
$ws.Range("K99").Value = 1976.3334
$ws.Range("L99").Value = 3605.2856
$ws.Range("M99").Value = -478.3334
$ws.Range("N99").Value = -6601.2856

# Row 107: The Gold Experience
$ws.Range("H107").Value = 2558.4473
$ws.Range("I107").Value = 1904.3572
$ws.Range("J107").Value = 4389.9
$ws.Range("K107").Value = 1904.3572
$ws.Range("L107").Value = 4389.9
$ws.Range("M107").Value = 15.64280000000008
$ws.Range("N107").Value = -8229.9

# Row 141: Awl Dreams Come True
$ws.Range("H141").Value = 36589.8
$ws.Range("J141").Value = 27649.666
$ws.Range("L141").Value = 27649.666
$ws.Range("N141").Value = -38009.666

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1262.5
$ws.Range("I22").Value = 330.2
$ws.Range("J22").Value = 1928.4286
$ws.Range("K22").Value = 330.2
$ws.Range("L22").Value = 1928.4286
$ws.Range("M22").Value = 19.80000000000001
$ws.Range("N22").Value = -2628.4286

# Row 38: Knock on Wood
$ws.Range("H38").Value = 500
$ws.Range("I38").Value = 500
$ws.Range("K38").Value = 500
$ws.Range("M38").Value = -123

# Row 46: Flintstone Fight
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 500
$ws.Range("M46").Value = -289

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1217.2858
$ws.Range("I132").Value = 879.47455
$ws.Range("J132").Value = 6200
$ws.Range("K132").Value = 2638.42365
$ws.Range("L132").Value = 18600
$ws.Range("M132").Value = -108.4236500000002
$ws.Range("N132").Value = -23660

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1557.1522
$ws.Range("I134").Value = 902.875
$ws.Range("J134").Value = 5919
$ws.Range("K134").Value = 2708.625
$ws.Range("L134").Value = 17757
$ws.Range("M134").Value = -173.625
$ws.Range("N134").Value = -22827

$ws = $wb.Worksheets.Item("CUL")
# Row 32: Convalescence Precedes Essence
$ws.Range("H32").Value = 1647.7
$ws.Range("I32").Value = 300
$ws.Range("J32").Value = 1984.625
$ws.Range("K32").Value = 900
$ws.Range("L32").Value = 5953.875
$ws.Range("M32").Value = -617
$ws.Range("N32").Value = -6519.875

# Row 87: Soup That Eats Like a Knight
$ws.Range("H87").Value = 5852.857
$ws.Range("I87").Value = 2207.3333
$ws.Range("J87").Value = 14966.667
$ws.Range("K87").Value = 6621.999899999999
$ws.Range("L87").Value = 44900.001
$ws.Range("M87").Value = -5373.999899999999
$ws.Range("N87").Value = -47396.001

# Row 90: Like Ma Used to Make (L)
$ws.Range("H90").Value = 5852.857
$ws.Range("I90").Value = 2207.3333
$ws.Range("J90").Value = 14966.667
$ws.Range("K90").Value = 19865.9997
$ws.Range("L90").Value = 134700.003
$ws.Range("M90").Value = -13625.9997
$ws.Range("N90").Value = -147180.003

# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 1991.4286
$ws.Range("I136").Value = 1408
$ws.Range("J136").Value = 3450
$ws.Range("K136").Value = 4224
$ws.Range("L136").Value = 10350
$ws.Range("M136").Value = 876
$ws.Range("N136").Value = -20550

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 918.72
$ws.Range("I107").Value = 429.2
$ws.Range("J107").Value = 1245.0667
$ws.Range("K107").Value = 429.2
$ws.Range("L107").Value = 1245.0667
$ws.Range("M107").Value = 1490.8
$ws.Range("N107").Value = -5085.0667

# Row 108: Satisfactory Sewing
$ws.Range("H108").Value = 28500
$ws.Range("J108").Value = 28500
$ws.Range("L108").Value = 28500
$ws.Range("N108").Value = -36180

# Row 132: On Board for Lar
$ws.Range("H132").Value = 1636.3442
$ws.Range("I132").Value = 1386.4222
$ws.Range("J132").Value = 2339.25
$ws.Range("K132").Value = 4159.2666
$ws.Range("L132").Value = 7017.75
$ws.Range("M132").Value = -1629.2666
$ws.Range("N132").Value = -12077.75

# Row 135: Fan of the Foreign
$ws.Range("H135").Value = 29500
$ws.Range("J135").Value = 29500
$ws.Range("L135").Value = 29500
$ws.Range("N135").Value = -39640

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 1333.1666
$ws.Range("I7").Value = 1319.8
$ws.Range("J7").Value = 1400
$ws.Range("K7").Value = 1319.8
$ws.Range("L7").Value = 1400
$ws.Range("M7").Value = -1207.8
$ws.Range("N7").Value = -1624

# Row 16: Saddle Sore
$ws.Range("H16").Value = 200001820
$ws.Range("I16").Value = 200001820
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 200001820
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -200001650
$ws.Range("N16").ClearContents()

# Row 126: Battered Books
$ws.Range("H126").Value = 1333.1666
$ws.Range("I126").Value = 1319.8
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 3959.4
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -1489.4
$ws.Range("N126").Value = -9140

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 2317.742
$ws.Range("I132").Value = 1558.6957
$ws.Range("K132").Value = 4676.0871
$ws.Range("M132").Value = -2146.0871

# Row 141: Just Generally Freezing
$ws.Range("H141").Value = 29933.334
$ws.Range("J141").Value = 29933.334
$ws.Range("L141").Value = 29933.334
$ws.Range("N141").Value = -40293.334

$ws = $wb.Worksheets.Item("WVR")
# Row 92: Modest Beginnings
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 107: Flax Wax
$ws.Range("H107").Value = 604.28125
$ws.Range("I107").Value = 290.68
$ws.Range("J107").Value = 1724.2858
$ws.Range("K107").Value = 872.04
$ws.Range("L107").Value = 5172.857400000001
$ws.Range("M107").Value = 1047.96
$ws.Range("N107").Value = -9012.857400000001

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1041.0233
$ws.Range("I136").Value = 507.77777
$ws.Range("J136").Value = 1940.875
$ws.Range("K136").Value = 1523.33331
$ws.Range("L136").Value = 5822.625
$ws.Range("M136").Value = 1026.66669
$ws.Range("N136").Value = -10922.625
